$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ID do Grupo"
$ws.Range("B1").Value = "Nome do Grupo"
$ws.Range("C1").Value = "QTD de Regras"
$ws.Range("D1").Value = "QTD de Recursos"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Grupo 1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Grupo 2"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0

# Build the header look (bold, centered, thin boxed border) on a scratch
# cell first, then fan it out to the header row in one paste so the whole
# range picks up a single finished style instead of several transient ones.
$tmp = $ws.Range("F1")
$tmp.Value = "x"
$tmp.Font.Bold = $true
$tmp.HorizontalAlignment = -4108
$tmp.VerticalAlignment = -4160
$tmp.Borders.LineStyle = 1

$tmp.Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$tmp.Clear()

# Column widths (best fit to content, like Excel's AutoFit would produce)
$ws.Columns.Item(1).ColumnWidth = 11.109375
$ws.Columns.Item(2).ColumnWidth = 14.5546875
$ws.Columns.Item(3).ColumnWidth = 13.21875
$ws.Columns.Item(4).ColumnWidth = 15.109375


# Match the saved selection
[void]$ws.Range("A3").Select()
